$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest cryptos snapshot.
# Column D values are stored as text (not numbers), so we temporarily force a
# text number-format before assigning, then clear formatting so the cell keeps
# its original (unstyled) appearance while retaining the text value.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.610.51'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.570.96'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.57'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.77'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.035.42'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.537.06'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.571.55'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.64'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.65'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.24'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('E25').Value = '  -4.03%  '
$ws.Range('E26').Value = '  +1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.93'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('E29').Value = '  -3.37%  '
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0802'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '454.52'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.24'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.85'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('E40').Value = '  -3.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '158.58'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.627'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.78'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0533'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0958'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('E47').Value = '  -3.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.96'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.42'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.959'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.58%  '
